$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 21.81010410189272
$ws.Range("C2").Value = 14.24849248060797
$ws.Range("D2").Value = 5.509729791891682
$ws.Range("E2").Value = 10.98168900857511
$ws.Range("F2").Value = 48.70175862158809
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.836636151400951
$ws.Range("M2").Value = 20.67638889649252
$ws.Range("N2").Value = 20.99764970910271

# Row 3
$ws.Range("B3").Value = 21.40109903940496
$ws.Range("C3").Value = 13.87932588755973
$ws.Range("D3").Value = 5.503787484904474
$ws.Range("E3").Value = 11.00212571772387
$ws.Range("F3").Value = 48.48214296658387
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.858482385446681
$ws.Range("M3").Value = 20.57850992337493
$ws.Range("N3").Value = 21.05364744127042

# Row 4
$ws.Range("B4").Value = 21.15314531910602
$ws.Range("C4").Value = 13.65258548663741
$ws.Range("D4").Value = 5.500147929783818
$ws.Range("E4").Value = 11.01581593387948
$ws.Range("F4").Value = 48.35987857803795
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.872916111439675
$ws.Range("M4").Value = 20.52342341889544
$ws.Range("N4").Value = 21.09001153987289

# Row 5
$ws.Range("B5").Value = 21.05304857057327
$ws.Range("C5").Value = 13.56032451679647
$ws.Range("D5").Value = 5.498666784904231
$ws.Range("E5").Value = 11.02168240543774
$ws.Range("F5").Value = 48.3132430087334
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 9.879054763923893
$ws.Range("M5").Value = 20.50225097985881
$ws.Range("N5").Value = 21.10532821433293

# Row 6
$ws.Range("B6").Value = 21.03648904635406
$ws.Range("C6").Value = 13.54501760348462
$ws.Range("D6").Value = 5.498420955065243
$ws.Range("E6").Value = 11.02267391145042
$ws.Range("F6").Value = 48.30569239407041
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 9.880089600751681
$ws.Range("M6").Value = 20.49881278380541
$ws.Range("N6").Value = 21.10790160975601

# Row 7
$ws.Range("B7").Value = 21.15179135265497
$ws.Range("C7").Value = 13.65134045673052
$ws.Range("D7").Value = 5.500127946790833
$ws.Range("E7").Value = 11.01589388611142
$ws.Range("F7").Value = 48.35923669647582
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.872997859312866
$ws.Range("M7").Value = 20.52313269524832
$ws.Range("N7").Value = 21.09021608969863

# Row 8
$ws.Range("B8").Value = 21.66850823854836
$ws.Range("C8").Value = 14.12130996993027
$ws.Range("D8").Value = 5.507678841596474
$ws.Range("E8").Value = 10.98849882562413
$ws.Range("F8").Value = 48.6234402689394
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.843957196407263
$ws.Range("M8").Value = 20.64160995037734
$ws.Range("N8").Value = 21.01654639423714

# Row 9
$ws.Range("B9").Value = 22.70054779713384
$ws.Range("C9").Value = 15.03584260755807
$ws.Range("D9").Value = 5.522570844890966
$ws.Range("E9").Value = 10.94381942191714
$ws.Range("F9").Value = 49.24006855333496
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.795089323030574
$ws.Range("M9").Value = 20.91294263563611
$ws.Range("N9").Value = 20.88780921635381

# Row 10
$ws.Range("B10").Value = 23.46174464826234
$ws.Range("C10").Value = 15.69523779890631
$ws.Range("D10").Value = 5.533582620902671
$ws.Range("E10").Value = 10.91648042150866
$ws.Range("F10").Value = 49.75114427578126
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.764095036284065
$ws.Range("M10").Value = 21.13492838551378
$ws.Range("N10").Value = 20.8028259483961

# Row 11
$ws.Range("B11").Value = 23.80695632005099
$ws.Range("C11").Value = 15.99093408569718
$ws.Range("D11").Value = 5.538611696734836
$ws.Range("E11").Value = 10.90522926210824
$ws.Range("F11").Value = 49.99572949702041
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.75105730993824
$ws.Range("M11").Value = 21.24055631831302
$ws.Range("N11").Value = 20.76625251031169

# Row 12
$ws.Range("B12").Value = 23.93739063236699
$ws.Range("C12").Value = 16.10217564245128
$ws.Range("D12").Value = 5.540519244546317
$ws.Range("E12").Value = 10.90113877868098
$ws.Range("F12").Value = 50.09003683200211
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.746272678379215
$ws.Range("M12").Value = 21.28119758785803
$ws.Range("N12").Value = 20.75270346894305

# Row 13
$ws.Range("B13").Value = 23.9093142770461
$ws.Range("C13").Value = 16.07825213671297
$ws.Range("D13").Value = 5.540108275200683
$ws.Range("E13").Value = 10.90201217898588
$ws.Range("F13").Value = 50.06965182347029
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.747296356306641
$ws.Range("M13").Value = 21.27241662114581
$ws.Range("N13").Value = 20.75560812173829

# Row 14
$ws.Range("B14").Value = 23.81769380172192
$ws.Range("C14").Value = 16.00010131612775
$ws.Range("D14").Value = 5.538768567034034
$ws.Range("E14").Value = 10.90488932900062
$ws.Range("F14").Value = 50.0034546877974
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.750660620787368
$ws.Range("M14").Value = 21.24388717986004
$ws.Range("N14").Value = 20.7651317960729

# Row 15
$ws.Range("B15").Value = 23.76153190371534
$ws.Range("C15").Value = 15.95213296621741
$ws.Range("D15").Value = 5.537948375891625
$ws.Range("E15").Value = 10.90667380528013
$ws.Range("F15").Value = 49.96312531609577
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.752741180412405
$ws.Range("M15").Value = 21.22649493015234
$ws.Range("N15").Value = 20.77100447684615

# Row 16
$ws.Range("B16").Value = 23.43915102273006
$ws.Range("C16").Value = 15.67581713666253
$ws.Range("D16").Value = 5.533254382552458
$ws.Range("E16").Value = 10.91723953338468
$ws.Range("F16").Value = 49.73539915179212
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.764968429591725
$ws.Range("M16").Value = 21.12811665751051
$ws.Range("N16").Value = 20.80525813029793

# Row 17
$ws.Range("B17").Value = 23.24100938978483
$ws.Range("C17").Value = 15.50512615645318
$ws.Range("D17").Value = 5.53038004656355
$ws.Range("E17").Value = 10.92402462363378
$ws.Range("F17").Value = 49.59875906442988
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.772741239071797
$ws.Range("M17").Value = 21.06893666852652
$ws.Range("N17").Value = 20.82680627882781

# Row 18
$ws.Range("B18").Value = 23.12695156198499
$ws.Range("C18").Value = 15.40655418597589
$ws.Range("D18").Value = 5.52872864633881
$ws.Range("E18").Value = 10.9280388350661
$ws.Range("F18").Value = 49.52130895741067
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.777311900359711
$ws.Range("M18").Value = 21.03533678495349
$ws.Range("N18").Value = 20.83939649931422

# Row 19
$ws.Range("B19").Value = 23.0883220307389
$ws.Range("C19").Value = 15.37311534036566
$ws.Range("D19").Value = 5.528169817750749
$ws.Range("E19").Value = 10.92941715840659
$ws.Range("F19").Value = 49.49528325618235
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.778876619978702
$ws.Range("M19").Value = 21.02403658434157
$ws.Range("N19").Value = 20.84369302854667

# Row 20
$ws.Range("B20").Value = 23.26211244231243
$ws.Range("C20").Value = 15.52333827820304
$ws.Range("D20").Value = 5.530685831726959
$ws.Range("E20").Value = 10.9232907910864
$ws.Range("F20").Value = 49.61318683560546
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.771903467606309
$ws.Range("M20").Value = 21.07519123957959
$ws.Range("N20").Value = 20.82449212175531

# Row 21
$ws.Range("B21").Value = 23.84461388456661
$ws.Range("C21").Value = 16.02307687725687
$ws.Range("D21").Value = 5.539161983975621
$ws.Range("E21").Value = 10.90403962727826
$ws.Range("F21").Value = 50.02285299637189
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.749668318287478
$ws.Range("M21").Value = 21.25224974863941
$ws.Range("N21").Value = 20.76232630103934

# Row 22
$ws.Range("B22").Value = 24.22356383586934
$ws.Range("C22").Value = 16.34536877018077
$ws.Range("D22").Value = 5.544720120600739
$ws.Range("E22").Value = 10.89244909323283
$ws.Range("F22").Value = 50.30040776146078
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.736024999040371
$ws.Range("M22").Value = 21.37169856762194
$ws.Range("N22").Value = 20.72344933794138

# Row 23
$ws.Range("B23").Value = 24.02151473150101
$ws.Range("C23").Value = 16.17378720003482
$ws.Range("D23").Value = 5.541751845905027
$ws.Range("E23").Value = 10.89854461082573
$ws.Range("F23").Value = 50.15139125677268
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.743225450514403
$ws.Range("M23").Value = 21.30761392162664
$ws.Range("N23").Value = 20.74403819446995

# Row 24
$ws.Range("B24").Value = 23.252572188941
$ws.Range("C24").Value = 15.51510593763762
$ws.Range("D24").Value = 5.530547582870415
$ws.Range("E24").Value = 10.92362220362124
$ws.Range("F24").Value = 49.60666059118471
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.772281906457913
$ws.Range("M24").Value = 21.0723622270214
$ws.Range("N24").Value = 20.82553772332906

# Row 25
$ws.Range("B25").Value = 22.4202342129273
$ws.Range("C25").Value = 14.7900735140774
$ws.Range("D25").Value = 5.518531146329877
$ws.Range("E25").Value = 10.95494097650194
$ws.Range("F25").Value = 49.06289827814515
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.807446087105982
$ws.Range("M25").Value = 20.8354758208491
$ws.Range("N25").Value = 20.92095114036796
